# "right display of orders"
# The due_date column (E) was showing a stray text date ("07.03,2020")
# and Excel date-serials for the other rows; normalize it to the plain
# number 4 for every order row. Also fix rows 4/5 which were incorrectly
# marked as finished ("done"/TRUE) - they are still "onTime"/FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the due_date column (E2:E5) to plain numbers and drop the
# leftover date number formatting so the cells show as plain numbers.
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("E2:E5").Style = "Normal"

# Rows 4 and 5 are not actually done yet - correct the status/done flag.
$ws.Range("H4").Value = $false
$ws.Range("I4").Value = "onTime"

$ws.Range("H5").Value = $false
$ws.Range("I5").Value = "onTime"

# Restore the user's last selection.
$ws.Range("J5").Select() | Out-Null
